$d = $word.ActiveDocument

# 1) "Lecture 05" -> "Lecture 07" (title heading)
$d.Content.Find.Execute("Lecture 05", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Lecture 07", 2)

# 2) " as lecture 4" -> " as lecture 6"
$d.Content.Find.Execute(" as lecture 4", $false, $false, $false, $false, $false, `
                         $true, 1, $false, " as lecture 6", 2)

# 3) Move the auto-managed "_GoBack" bookmark (last-edit-position marker) from
#    the end of the document (around the final inserted image) to right after
#    the text that was just edited ("... as lecture 6"), matching where Word
#    would naturally drop it after the last text edit made to the document.
$find = $d.Content
$find.Find.Execute(" as lecture 6", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$goBackRange = $d.Range($find.End, $find.End)
$d.Bookmarks.Add("_GoBack", $goBackRange)
